# Finished Week 13 logging
# Update the "H" (home) row totals on both the OFF and DEF sheets with the
# latest weekly stats: Short Att, Short Comp, Deep Att, Deep Comp, Short Int.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 ("H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 148   # Short Att
$wsOff.Range("C2").Value = 100   # Short Comp
$wsOff.Range("D2").Value = 36    # Deep Att
$wsOff.Range("F2").Value = 4     # Short Int

# --- DEF sheet: row 2 ("H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 125   # Short Att
$wsDef.Range("C2").Value = 84    # Short Comp
$wsDef.Range("D2").Value = 23    # Deep Att
$wsDef.Range("E2").Value = 17    # Deep Comp

$wb.Save()
